# Adds a new "affix_type" column (AN) to the Affixes sheet, filled with
# the value 4 for every existing data row, to support affix-type filters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in AN1
$ws.Range("AN1").Value = "affix_type"

# New data values for AN2:AN25 (all existing affix rows)
$ws.Range("AN2:AN25").Value = 4

# Match the new column's width to the other data columns
$ws.Columns.Item(40).ColumnWidth = 16.3333333

# Update the active selection to the newly added cells, like Excel would
# leave it after filling down a new column
$ws.Range("AN22:AN25").Select()
